$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RegressionTests")
$ws.Range("A1").Value = "test"
